$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The s2cDNADate columns (A = s1cDNADate, D = s2cDNADate) for rows 24-45
# change from "01.09.17" to "01.09.18". Force text (not date) entry by
# temporarily marking the range as Text before assignment, then restore
# the default "Normal" style so the cell keeps its original (unstyled)
# appearance while remaining a text value. (Applied to column A and D
# separately -- a multi-area Range only reliably applies formatting to
# its first area.)
$colA = $ws.Range("A24:A45")
$colD = $ws.Range("D24:D45")
$colA.NumberFormat = "@"
$colD.NumberFormat = "@"
for ($r = 24; $r -le 45; $r++) {
    $ws.Range("A$r").Value = "01.09.18"
    $ws.Range("D$r").Value = "01.09.18"
}
$colA.Style = "Normal"
$colD.Style = "Normal"

# Restore the view/selection state: scrolled down with the active cell at
# D24 and the selection spanning the newly-edited rows.
$excel.ActiveWindow.ScrollRow = 18
$ws.Range("D24:D45").Select()
